$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 9, shifting existing rows 9..114 down to 10..115.
$ws.Rows("9:9").Insert()

# Populate the newly inserted row 9 with the new record.
$ws.Range("A9").Value = 5
$ws.Range("B9").Value = "Macroferia Regional de Talca"
$ws.Range("C9").Value = "Maule"
$ws.Range("D9").Value = 44881
$ws.Range("E9").Value = 7
$ws.Range("F9").Value = 100112022
$ws.Range("G9").Value = "Arveja Verde"
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 500
$ws.Range("K9").Value = 17000
$ws.Range("L9").Value = 17000
$ws.Range("M9").Value = 17000
$ws.Range("N9").Value = "`$/saco 25 kilos"
$ws.Range("O9").Value = "Región del Maule"
$ws.Range("P9").Value = 680
$ws.Range("Q9").Value = 25
$ws.Range("R9").Value = "Hortaliza"
